$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1650.25
$ws.Range("J19").Value = 1533.3334
$ws.Range("L19").Value = 1533.3334
$ws.Range("N19").Value = -1883.3334

$ws.Range("H55").Value = 1421.5
$ws.Range("I55").Value = 975.8
$ws.Range("J55").Value = 1867.2
$ws.Range("K55").Value = 975.8
$ws.Range("L55").Value = 1867.2
$ws.Range("M55").Value = -761.8
$ws.Range("N55").Value = -2295.2

$ws.Range("H96").Value = 834.5
$ws.Range("J96").Value = 1223.25
$ws.Range("L96").Value = 3669.75
$ws.Range("N96").Value = -6415.75

$ws.Range("H104").Value = 755.1429000000001
$ws.Range("I104").Value = 755.1429000000001
$ws.Range("K104").Value = 2265.4287
$ws.Range("M104").Value = -518.4287000000004

$ws.Range("H113").Value = 11131.25
$ws.Range("J113").Value = 14520
$ws.Range("L113").Value = 14520
$ws.Range("N113").Value = -21028

$ws.Range("H137").Value = 2138.9443
$ws.Range("I137").Value = 806.1
$ws.Range("J137").Value = 3805
$ws.Range("K137").Value = 2418.3
$ws.Range("L137").Value = 11415
$ws.Range("M137").Value = 131.6999999999998
$ws.Range("N137").Value = -16515

$ws.Range("H141").Value = 1806.3334
$ws.Range("I141").Value = 1361
$ws.Range("J141").Value = 3365
$ws.Range("K141").Value = 4083
$ws.Range("L141").Value = 10095
$ws.Range("M141").Value = 1097
$ws.Range("N141").Value = -20455

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H88").Value = 1396.8462
$ws.Range("I88").Value = 1285.75
$ws.Range("J88").Value = 1574.6
$ws.Range("K88").Value = 1285.75
$ws.Range("L88").Value = 1574.6
$ws.Range("M88").Value = -879.75
$ws.Range("N88").Value = -2386.6

$ws.Range("H91").Value = 1396.8462
$ws.Range("I91").Value = 1285.75
$ws.Range("J91").Value = 1574.6
$ws.Range("K91").Value = 1285.75
$ws.Range("L91").Value = 1574.6
$ws.Range("M91").Value = 118.25
$ws.Range("N91").Value = -4382.6

$ws.Range("H97").Value = 756.3
$ws.Range("I97").Value = 729.2222
$ws.Range("K97").Value = 729.2222
$ws.Range("M97").Value = -233.2222

$ws.Range("H122").Value = 2992
$ws.Range("I122").Value = 2995
$ws.Range("K122").Value = 8985
$ws.Range("M122").Value = -6535

$ws.Range("H132").Value = 3050
$ws.Range("I132").Value = 3079
$ws.Range("J132").Value = 2977.5
$ws.Range("K132").Value = 9237
$ws.Range("L132").Value = 8932.5
$ws.Range("M132").Value = -6707
$ws.Range("N132").Value = -13992.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2247
$ws.Range("I99").Value = 1896.4
$ws.Range("K99").Value = 1896.4
$ws.Range("M99").Value = -398.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 54500
$ws.Range("J97").Value = 54500
$ws.Range("L97").Value = 54500
$ws.Range("N97").Value = -56482

$ws.Range("H132").Value = 2013.3636
$ws.Range("I132").Value = 2013.3636
$ws.Range("K132").Value = 6040.0908
$ws.Range("M132").Value = -3510.0908

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()

$ws.Range("H46").Value = 9333
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 9333
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 27999
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -28181

$ws.Range("H87").Value = 13250
$ws.Range("I87").Value = 13250
$ws.Range("K87").Value = 39750
$ws.Range("M87").Value = -38502

$ws.Range("H90").Value = 13250
$ws.Range("I90").Value = 13250
$ws.Range("K90").Value = 119250
$ws.Range("M90").Value = -113010

$ws.Range("H114").Value = 1083.091
$ws.Range("I114").Value = 1443.5
$ws.Range("J114").Value = 877.1429000000001
$ws.Range("K114").Value = 4330.5
$ws.Range("L114").Value = 2631.4287
$ws.Range("M114").Value = -1076.5
$ws.Range("N114").Value = -9139.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2000
$ws.Range("J70").Value = 2000
$ws.Range("L70").Value = 2000
$ws.Range("N70").Value = -2540

$ws.Range("H73").Value = 2000
$ws.Range("J73").Value = 2000
$ws.Range("L73").Value = 2000
$ws.Range("N73").Value = -3872

$ws.Range("H96").Value = 14000
$ws.Range("J96").Value = 14000
$ws.Range("L96").Value = 14000
$ws.Range("N96").Value = -19492

$ws.Range("H97").Value = 682.2308
$ws.Range("I97").Value = 496.125
$ws.Range("J97").Value = 980
$ws.Range("K97").Value = 496.125
$ws.Range("L97").Value = 980
$ws.Range("M97").Value = -0.125
$ws.Range("N97").Value = -1972

$ws.Range("H122").Value = 657.75
$ws.Range("I122").Value = 477
$ws.Range("K122").Value = 1431
$ws.Range("M122").Value = 1019

$ws.Range("H126").Value = 8613.857
$ws.Range("I126").Value = 8099.4
$ws.Range("K126").Value = 24298.2
$ws.Range("M126").Value = -21828.2

$ws.Range("H132").Value = 1887.1428
$ws.Range("I132").Value = 1859.2307
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 5577.6921
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -3047.6921
$ws.Range("N132").Value = -11810

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4385.2
$ws.Range("J82").Value = 5306.5
$ws.Range("L82").Value = 5306.5
$ws.Range("N82").Value = -6028.5

$ws.Range("H85").Value = 4385.2
$ws.Range("J85").Value = 5306.5
$ws.Range("L85").Value = 5306.5
$ws.Range("N85").Value = -7802.5

$ws.Range("H93").Value = 940.58826
$ws.Range("I93").Value = 913.6429000000001
$ws.Range("K93").Value = 913.6429000000001
$ws.Range("M93").Value = 334.3570999999999
